$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block (rows 2-15): Technology, NumCores, Runtime
$data = @(
    @("serial", 1, 880),
    @("serial", 2, 880),
    @("serial", 4, 880),
    @("serial", 8, 880),
    @("serial", 12, 880),
    @("open_mp", 1, 1443),
    @("open_mp", 2, 1625),
    @("open_mp", 4, 1646),
    @("open_mp", 8, 810),
    @("open_mp", 12, 609),
    @("mpi", 2, 1071),
    @("mpi", 4, 727),
    @("mpi", 8, 650),
    @("mpi", 12, 653)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
